$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H header "Save", formatted like the rest of the header row (B1:G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# "Save" flag values for rows 2-17 (0/1)
$saveValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 1
    15 = 1
    16 = 0
    17 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
